$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, bordered) from an existing header cell to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Add new header labels for team record
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in team record values for each data row (2 through 48)
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 72  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 90  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
